# Weekly update: insert a new Puerro (leek) price-survey record for
# "Vega Modelo de Temuco" as row 224, pushing the existing rows 224-241
# down to 225-242 (dimension grows from A1:R241 to A1:R242).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 224..241 down by one, leaving a blank row 224 to populate.
$ws.Rows("224:224").Insert()

# Populate the new row 224 with the latest week's observation.
$ws.Cells.Item(224, 1).Value = 10
$ws.Cells.Item(224, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(224, 3).Value = "La Araucanía"
$ws.Cells.Item(224, 4).Value = 44826
$ws.Cells.Item(224, 5).Value = 9
$ws.Cells.Item(224, 6).Value = 100112005
$ws.Cells.Item(224, 7).Value = "Puerro"
$ws.Cells.Item(224, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 20
$ws.Cells.Item(224, 11).Value = 15000
$ws.Cells.Item(224, 12).Value = 15000
$ws.Cells.Item(224, 13).Value = 15000
$ws.Cells.Item(224, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(224, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(224, 16).Value = 1250
$ws.Cells.Item(224, 17).Value = 12
$ws.Cells.Item(224, 18).Value = "Hortaliza"
